# Append two new daily data rows (2025-10-26) for both stations to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting of the last two existing rows (50:51) down onto the
# new rows (52:53) so the new cells reuse the same cell styles instead of
# generating new number-format entries.
$ws.Range("A50:F51").Copy()
$ws.Range("A52:F53").PasteSpecial(-4122) # xlPasteFormats

# Row 52: 四方坪站 (station in column B already uses shared string index 4)
$ws.Cells.Item(52, 1).Value = 45956
$ws.Cells.Item(52, 2).Value = "四方坪站"
$ws.Cells.Item(52, 3).Value = 9695.2000000000007
$ws.Cells.Item(52, 4).Value = 8071.5
$ws.Cells.Item(52, 5).Value = 3347.98
$ws.Cells.Item(52, 6).Value = 409

# Row 53: 高岭站
$ws.Cells.Item(53, 1).Value = 45956
$ws.Cells.Item(53, 2).Value = "高岭站"
$ws.Cells.Item(53, 3).Value = 4219.62
$ws.Cells.Item(53, 4).Value = 3334.28
$ws.Cells.Item(53, 5).Value = 1092.1600000000001
$ws.Cells.Item(53, 6).Value = 146

# Update the selection/active cell to match the saved view state
$ws.Range("H52").Select()
